# "April 3 - Maya" commit: airplane data spreadsheet completed with 10
# planes total (1 header row + 10 data rows). Economy "total seats" header
# is renamed to be explicit, and 10 aircraft/airline rows are added with
# their seating data across economy / business / first class sections.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
# Only real content change here: "total seats" -> "total seats - economy"
$ws.Range("A1").Value = "aircraft_type"
$ws.Range("B1").Value = "airline"
$ws.Range("C1").Value = "total seats - economy"
$ws.Range("D1").Value = "seat width - economy"
$ws.Range("E1").Value = "seat depth - economy"
$ws.Range("F1").Value = "seat arrangement - economy"
$ws.Range("G1").Value = "total seats - economy"
$ws.Range("H1").Value = "seat width - business"
$ws.Range("I1").Value = "seat depth - business"
$ws.Range("J1").Value = "seat arrangement - business"
$ws.Range("K1").Value = "total seats - business"
$ws.Range("L1").Value = "seat width - first class"
$ws.Range("M1").Value = "seat depth - first class"
$ws.Range("N1").Value = "seat arrangement - first class"
$ws.Range("O1").Value = "total seats - first class"

# --- Existing data rows (2-3), unchanged content but re-asserted ------
$ws.Range("A2").Value = "A220-300"
$ws.Range("B2").Value = "Air Canada"
$ws.Range("C2").Value = 137
$ws.Range("D2").Value = 19
$ws.Range("E2").Value = 30
$ws.Range("F2").Value = "2 3"
$ws.Range("G2").Value = 125
$ws.Range("H2").Value = 21
$ws.Range("I2").Value = 37
$ws.Range("J2").Value = "2 2"
$ws.Range("K2").Value = 12
$ws.Range("O2").Value = 0

$ws.Range("A3").Value = "Bombardier CRJ-900"
$ws.Range("B3").Value = "Air Canada"
$ws.Range("C3").Value = 76
$ws.Range("D3").Value = 17
$ws.Range("E3").Value = 31
$ws.Range("F3").Value = "2 2"
$ws.Range("G3").Value = 64
$ws.Range("H3").Value = 21
$ws.Range("I3").Value = 37
$ws.Range("J3").Value = "1 2"
$ws.Range("K3").Value = 12
$ws.Range("O3").Value = 0

# --- New aircraft rows (4-11) ------------------------------------------
$ws.Range("A4").Value = "Beechcraft 1900D"
$ws.Range("B4").Value = "Air Canada"
$ws.Range("C4").Value = 18
$ws.Range("D4").Value = 20.2
$ws.Range("E4").Value = 30
$ws.Range("F4").Value = "1 1"
$ws.Range("G4").Value = 18
$ws.Range("K4").Value = 0
$ws.Range("O4").Value = 0

$ws.Range("A5").Value = "Boeing 737-300 (733)"
$ws.Range("B5").Value = "Lion Airlines"
$ws.Range("C5").Value = 150
$ws.Range("D5").Value = 17
$ws.Range("E5").Value = 30
$ws.Range("F5").Value = "3 3"
$ws.Range("G5").Value = 150
$ws.Range("K5").Value = 0
$ws.Range("O5").Value = 0

$ws.Range("A6").Value = "Airbus A319 (319) Layout 1"
$ws.Range("B6").Value = "United"
$ws.Range("C6").Value = 128
$ws.Range("D6").Value = 17.7
$ws.Range("E6").Value = 30
$ws.Range("F6").Value = "3 3"
$ws.Range("G6").Value = 78
$ws.Range("H6").Value = 17.7
$ws.Range("I6").Value = 34
$ws.Range("J6").Value = "3 3"
$ws.Range("K6").Value = 42
$ws.Range("L6").Value = 20.5
$ws.Range("M6").Value = 37
$ws.Range("N6").Value = "2 2"
$ws.Range("O6").Value = 8

$ws.Range("A7").Value = "Embraer E-175 (E75) Layout 1"
$ws.Range("B7").Value = "Delta"
$ws.Range("C7").Value = 76
$ws.Range("D7").Value = 18.25
$ws.Range("E7").Value = 31
$ws.Range("F7").Value = "2 2"
$ws.Range("G7").Value = 52
$ws.Range("H7").Value = 18.25
$ws.Range("I7").Value = 34
$ws.Range("J7").Value = "2 2"
$ws.Range("K7").Value = 12
$ws.Range("L7").Value = 20
$ws.Range("M7").Value = 37
$ws.Range("N7").Value = "1 2"
$ws.Range("O7").Value = 12

$ws.Range("A8").Value = "Airbus A320-200 (320)"
$ws.Range("B8").Value = "Condor"
$ws.Range("C8").Value = 164
$ws.Range("D8").Value = 17
$ws.Range("E8").Value = 29
$ws.Range("F8").Value = "3 3"
$ws.Range("G8").Value = 132
$ws.Range("H8").Value = 17
$ws.Range("I8").Value = 29
$ws.Range("J8").Value = "3 3"
$ws.Range("K8").Value = 32
$ws.Range("O8").Value = 0

$ws.Range("A9").Value = "Boeing 737-700 Layout 2"
$ws.Range("B9").Value = "SpiceJet"
$ws.Range("C9").Value = 134
$ws.Range("D9").Value = 17
$ws.Range("E9").Value = 28
$ws.Range("F9").Value = "3 3"
$ws.Range("G9").Value = 126
$ws.Range("H9").Value = 19
$ws.Range("I9").Value = 32
$ws.Range("J9").Value = "2 2"
$ws.Range("K9").Value = 8
$ws.Range("O9").Value = 0

$ws.Range("A10").Value = "ATR 72-600 (ATR)"
$ws.Range("B10").Value = "Air France"
$ws.Range("C10").Value = 72
$ws.Range("D10").Value = 17
$ws.Range("E10").Value = 31
$ws.Range("F10").Value = "2 2"
$ws.Range("G10").Value = 72
$ws.Range("K10").Value = 0
$ws.Range("O10").Value = 0

$ws.Range("A11").Value = "Boeing 717-200 (717)"
$ws.Range("B11").Value = "Hawaiian Airlines"
$ws.Range("C11").Value = 128
$ws.Range("D11").Value = 18
$ws.Range("E11").Value = 30
$ws.Range("F11").Value = "2 3"
$ws.Range("G11").Value = 120
$ws.Range("L11").Value = 18.5
$ws.Range("M11").Value = 37
$ws.Range("N11").Value = "2 2"
$ws.Range("O11").Value = 8

# --- Formatting ---------------------------------------------------------
# Wrap text on the long aircraft-name / header cells (matches existing
# wrap-text style already used for the "Bombardier CRJ-900" row).
$ws.Range("O1").WrapText = $true
$ws.Range("A7").WrapText = $true

# Row heights grow to fit the wrapped text.
$ws.Rows.Item(1).RowHeight = 29.25
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(7).RowHeight = 15

# --- Selection ------------------------------------------------------------
$ws.Range("O12").Select() | Out-Null
